{"js": "const replacements = [\n  [\"2024-11-15 Friday\", \"2024-11-16 Saturday\"],\n  [\"67\u00d737=\", \"20\u00d759=\"],\n  [\"57\u00d756=\", \"31\u00d745=\"],\n  [\"93\u00d769=\", \"83\u00d712=\"],\n  [\"56\u00d754=\", \"27\u00d738=\"],\n  [\"88\u00d721=\", \"81\u00d763=\"],\n  [\"67\u00d769=\", \"87\u00d738=\"],\n  [\"69\u00d716=\", \"82\u00d794=\"],\n  [\"38\u00d731=\", \"22\u00d740=\"],\n  [\"74\u00d785=\", \"57\u00d779=\"],\n  [\"87\u00d732=\", \"63\u00d762=\"],\n  [\"42\u00d767=\", \"85\u00d778=\"],\n  [\"88\u00d752=\", \"12\u00d738=\"],\n  [\"87\u00d714=\", \"36\u00d739=\"],\n  [\"48\u00d794=\", \"89\u00d749=\"],\n  [\"99\u00d716=\", \"99\u00d773=\"],\n  [\"26\u00d771=\", \"95\u00d725=\"],\n  [\"52\u00d794=\", \"46\u00d772=\"],\n  [\"19\u00d718=\", \"30\u00d764=\"],\n  [\"35\u00d771=\", \"82\u00d740=\"],\n  [\"74\u00d764=\", \"26\u00d776=\"],\n  [\"51\u00d759=\", \"22\u00d766=\"],\n  [\"28\u00d798=\", \"13\u00d733=\"],\n  [\"60\u00d795=\", \"91\u00d737=\"],\n  [\"32\u00d762=\", \"14\u00d721=\"],\n  [\"15\u00d711=\", \"59\u00d770=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nreturn \"done\";\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2024-11-15 Friday\"; Replace = \"2024-11-16 Saturday\" },\n    @{ Find = \"67\u00d737=\"; Replace = \"20\u00d759=\" },\n    @{ Find = \"57\u00d756=\"; Replace = \"31\u00d745=\" },\n    @{ Find = \"93\u00d769=\"; Replace = \"83\u00d712=\" },\n    @{ Find = \"56\u00d754=\"; Replace = \"27\u00d738=\" },\n    @{ Find = \"88\u00d721=\"; Replace = \"81\u00d763=\" },\n    @{ Find = \"67\u00d769=\"; Replace = \"87\u00d738=\" },\n    @{ Find = \"69\u00d716=\"; Replace = \"82\u00d794=\" },\n    @{ Find = \"38\u00d731=\"; Replace = \"22\u00d740=\" },\n    @{ Find = \"74\u00d785=\"; Replace = \"57\u00d779=\" },\n    @{ Find = \"87\u00d732=\"; Replace = \"63\u00d762=\" },\n    @{ Find = \"42\u00d767=\"; Replace = \"85\u00d778=\" },\n    @{ Find = \"88\u00d752=\"; Replace = \"12\u00d738=\" },\n    @{ Find = \"87\u00d714=\"; Replace = \"36\u00d739=\" },\n    @{ Find = \"48\u00d794=\"; Replace = \"89\u00d749=\" },\n    @{ Find = \"99\u00d716=\"; Replace = \"99\u00d773=\" },\n    @{ Find = \"26\u00d771=\"; Replace = \"95\u00d725=\" },\n    @{ Find = \"52\u00d794=\"; Replace = \"46\u00d772=\" },\n    @{ Find = \"19\u00d718=\"; Replace = \"30\u00d764=\" },\n    @{ Find = \"35\u00d771=\"; Replace = \"82\u00d740=\" },\n    @{ Find = \"74\u00d764=\"; Replace = \"26\u00d776=\" },\n    @{ Find = \"51\u00d759=\"; Replace = \"22\u00d766=\" },\n    @{ Find = \"28\u00d798=\"; Replace = \"13\u00d733=\" },\n    @{ Find = \"60\u00d795=\"; Replace = \"91\u00d737=\" },\n    @{ Find = \"32\u00d762=\"; Replace = \"14\u00d721=\" },\n    @{ Find = \"15\u00d711=\"; Replace = \"59\u00d770=\" }\n)\n\nforeach ($item in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($item.Find, $false, $false, $false, $false, $false, $true, 1, $false, $item.Replace, 2) | Out-Null\n}\n"}
